# Weekly refresh of the Fruta/Hortaliza (Frambuesa) price rows.
# The price-report blocks (2 rows each: Primera/Segunda quality) for
# rows 4-14 are reshuffled onto a new set of report dates, per the
# latest "Vega Monumental Concepción" weekly pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Primera) - date moves from 2021-12-29 to 2022-02-10; figures unchanged
$ws.Range("D4").Value = 44602

# Row 5 (Segunda) - same date move; figures unchanged
$ws.Range("D5").Value = 44602

# Row 6 (Primera) - date 2021-12-02 -> 2022-02-17, new price figures
$ws.Range("D6").Value = 44609
$ws.Range("N6").Value = 6500
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 6750
$ws.Range("S6").Value = 3375

# Row 7 (Segunda) - date 2021-12-02 -> 2022-02-17, new price figures
$ws.Range("D7").Value = 44609
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("S7").Value = 3000

# Row 8 (Primera) - date 2022-02-25 -> 2021-12-02, new price figures
$ws.Range("D8").Value = 44532
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 5000

# Row 9 - date 2022-02-10 -> 2021-12-02, quality Primera -> Segunda, new figures
$ws.Range("D9").Value = 44532
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 4000

# Row 10 - date 2022-02-10 -> 2020-12-30, quality Segunda -> Primera, new figures
$ws.Range("D10").Value = 44195
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 3000
$ws.Range("O10").Value = 3500
$ws.Range("P10").Value = 3250
$ws.Range("S10").Value = 1625

# Row 11 - date 2022-02-17 -> 2020-12-30, quality Primera -> Segunda, new figures
$ws.Range("D11").Value = 44195
$ws.Range("L11").Value = "Segunda"
$ws.Range("N11").Value = 2500
$ws.Range("O11").Value = 2500
$ws.Range("P11").Value = 2500
$ws.Range("S11").Value = 1250

# Row 12 - date 2022-02-17 -> 2022-02-25, quality Segunda -> Primera, new figures
$ws.Range("D12").Value = 44617
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6500
$ws.Range("S12").Value = 3250

# Row 13 (Primera) - date 2020-12-30 -> 2021-12-29, new price figures
$ws.Range("D13").Value = 44559
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 7000
$ws.Range("P13").Value = 6500
$ws.Range("S13").Value = 3250

# Row 14 (Segunda) - date 2020-12-30 -> 2021-12-29, new price figures
$ws.Range("D14").Value = 44559
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5000
$ws.Range("P14").Value = 5000
$ws.Range("S14").Value = 2500
